# Automatic update of PEBCOM sheet:
#  - Fills in the previously-empty OT value for case 6075 (row 63, col E)
#  - Appends a new case row (6137 / LA PLATA AV. 1058) as row 64

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text (matching how every other column
# A-L cell in this sheet is stored) without leaving a lingering
# "Text" number-format on the cell - Excel sets quotePrefix while the
# value is typed, ClearFormats() drops that again once the literal
# text value has been committed to the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 63: OT number had not been assigned yet - fill it in now.
Set-TextValue $ws.Range("E63") "807458159"

# Row 64: new case reported 6/12/2025.
Set-TextValue $ws.Range("A64") "6137"
Set-TextValue $ws.Range("B64") "6/12/2025"
Set-TextValue $ws.Range("C64") "LA PLATA AV. 1058"
Set-TextValue $ws.Range("D64") "7"
Set-TextValue $ws.Range("E64") "807458383"
Set-TextValue $ws.Range("F64") "PEBCOM"
Set-TextValue $ws.Range("G64") "Pendiente"
Set-TextValue $ws.Range("H64") "Ver con inspector tratar de colocar r400 para sacar las dos columnas terminales existentes evaluar en campo"
Set-TextValue $ws.Range("I64") "1"
Set-TextValue $ws.Range("J64") "Cambio"
Set-TextValue $ws.Range("K64") "Sin equipos"
Set-TextValue $ws.Range("L64") "Terminal"

$ws.Range("M64").Value = -58.426431
$ws.Range("N64").Value = -34.627954
